# Generate Report for Archive
#
# 1) The status string shown for the localized file moves from
#    "Ready for handoff" to "In Translation" everywhere it appears:
#    Overview!E2, Overview!F2, zh-cn!C2, de-de!C2.
# 2) Because the new label is shorter, the "Status" column is re-sized
#    (narrower) on all three sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- 1. Update the status text ---------------------------------------
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2. Re-size the Status columns to fit the new, shorter text ------
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
